# Registro de empresa (No modificar el excel)
#
# 1) Move the inline "Tipos de empresas / valor" drop-down list that lived in
#    columns M:N of "Registro de empresas" out to its own sheet named
#    "Lista desplegables".
# 2) Replace the in-sheet list with a proper list-based data validation on
#    column D (D2:D1048576) that points at the new sheet.
# 3) Tidy up the selection / view state that Excel re-writes as a side effect
#    of the edit.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Registro de empresas")

# ---------------------------------------------------------------------------
# 1. Remove the now-redundant "Tipos de empresas / valor" helper list that
#    lived in columns M:N of "Registro de empresas" (it's about to be
#    recreated on its own sheet).
# ---------------------------------------------------------------------------
$ws2.Range("M1:N21").Clear()

# ---------------------------------------------------------------------------
# 2. Create the new "Lista desplegables" sheet at the end of the workbook and
#    populate it with the values that used to sit in M1:N21 of sheet2.
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$listSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$listSheet.Name = "Lista desplegables"

# Copy the header formatting (bold + fill) from the existing header row so we
# reuse the same style instead of minting a new one, then set the values.
$ws2.Range("A1").Copy()
$listSheet.Range("A1:B1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$listSheet.Range("A1").Value = "Tipos de empresas"
$listSheet.Range("B1").Value = "valor"

$tipos = @(
    "Grupo de Investigación de Universidad",
    "Centro de I+D+i",
    "Desarrollo de software",
    "Fabricante de componentes",
    "Fabricante de módulos",
    "Fabricante de sistemas",
    "Ingeniería ",
    "Distribución de productos",
    "Consultoría de I+D+i",
    "startup",
    "Aceleradora",
    "Incubadora",
    "Venture capital",
    "Business Angel",
    "Corporate",
    "Empresa industrial usuaria de tecnología",
    "Hospital o centro sanitario",
    "Medio de comunicación",
    "Empresa de servicios",
    "Administración pública"
)

for ($i = 0; $i -lt $tipos.Count; $i++) {
    $row = $i + 2
    $listSheet.Cells.Item($row, 1).Value = $tipos[$i]
    $listSheet.Cells.Item($row, 2).Value = $i
}

$listSheet.Range("D2:D1048576").Select()

# ---------------------------------------------------------------------------
# 3. Add a list data validation on D2:D1048576 of "Registro de empresas" that
#    references the new sheet.
# ---------------------------------------------------------------------------
$validation = $ws2.Range("D2:D1048576").Validation
$validation.Add(3, 1, 1, "='Lista desplegables'!`$B`$2:`$B`$21")
$validation.ErrorTitle = "Entrada inválida"
$validation.ErrorMessage = "Seleccione un valor de la lista"
$validation.InputTitle = "Tipo de empresa"
$validation.InputMessage = "Seleccione un tipo de empresa"
$validation.ShowInput = $true
$validation.ShowError = $true

$ws2.Range("D3").Select()

# ---------------------------------------------------------------------------
# 4. Leave "Registro de empresas" as the active sheet (matches tabSelected in
#    the original file) and make sure the new sheet isn't left selected.
# ---------------------------------------------------------------------------
$ws2.Activate()
